$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Update the registration-ID prompt text (was: "The UUID of the client is:")
$ws.Range("G2").Value = "Please choose an ID for the patient:"

# Re-word remaining client-facing prompts to refer to "patient" instead of "client"
$ws.Range("G4").Value  = "What is your patient's first name?"
$ws.Range("G6").Value  = "What is your patient's middle name?"
$ws.Range("G9").Value  = "What is your patient's surname?"
$ws.Range("G12").Value = "What is your patient's assigned birth sex?"
$ws.Range("G15").Value = "When was your patient born?"
$ws.Range("G18").Value = "Where does your patient live?"
$ws.Range("G21").Value = "What country does your patient live in?"
$ws.Range("G24").Value = "What is your patient's telephone number?"
$ws.Range("G27").Value = "What is your patient's father's name?"
$ws.Range("G30").Value = "What is your patient's mother's name?"

# Update the active selection to reflect where editing ended up
$ws.Activate()
$ws.Range("G30").Select()
